$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "aliyunctf{smashi_your_SUIbscribe_button_now_do_it_quickly_pleaSEI_UC5CwaMl1eIgY8h02uZw7u8A}"
